$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.082.22'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.40%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.737.52'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '601.01'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.20%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '167.22'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.22%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.736.18'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E9').Value = '  +0.87%  '
$ws.Range('E10').Value = '  +3.67%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.37'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.11%  '
$ws.Range('E12').Value = '  +0.29%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '38.09'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.18%  '
$ws.Range('E14').Value = '  +1.30%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.362.54'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.01%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.713.33'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.65%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '69.022.03'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E18').Value = '  +1.46%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.35'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.15%  '
$ws.Range('E20').Value = '  -1.69%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.20'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +11.09%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '491.42'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.13%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.727'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.42%  '
$ws.Range('E24').Value = '  +7.11%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '84.74'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.19%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.29'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.93%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.28'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.67%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.05'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.30%  '
$ws.Range('E30').Value = '  +1.14%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.20'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.67%  '
$ws.Range('E32').Value = '  +0.44%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '31.42'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.77%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.882.83'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.01%  '
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.670.38'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.04%  '
$ws.Range('E37').Value = '  +0.00%  '
$ws.Range('B38').Value = 'Mantle'
$ws.Range('C38').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.01'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.03%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.139'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +5.42%  '
$ws.Range('E40').Value = '  +2.05%  '
$ws.Range('E41').Value = '  +0.02%  '
$ws.Range('E42').Value = '  +5.67%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '48.80'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.52%  '
$ws.Range('E44').Value = '  +0.23%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '422.96'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.46%  '
$ws.Range('E46').Value = '  +0.45%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '39.97'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.58%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '141.27'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.37%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.778.81'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.26%  '
$ws.Range('E51').Value = '  -0.07%  '
